$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1855072463768116
$ws.Range("C2").Value = 0.5710144927536231
$ws.Range("J2").Value = 0.008695652173913044
$ws.Range("P2").Value = 0.1246376811594203
$ws.Range("S2").Value = 0.1101449275362319
$ws.Range("B3").Value = 0.004878048780487805
$ws.Range("C3").Value = 0.02926829268292683
$ws.Range("J3").Value = 0.05365853658536585
$ws.Range("P3").Value = 0.7073170731707317
$ws.Range("S3").Value = 0.2048780487804878
$ws.Range("B6").Value = 0.1138613861386139
$ws.Range("D6").Value = 0.0198019801980198
$ws.Range("F6").Value = 0.03465346534653466
$ws.Range("J6").Value = 0.2722772277227723
$ws.Range("O6").Value = 0.02475247524752475
$ws.Range("Q6").Value = 0.1732673267326733
$ws.Range("R6").Value = 0.06930693069306931
$ws.Range("S6").Value = 0.2920792079207921
$ws.Range("B7").Value = 0.1685393258426966
$ws.Range("D7").Value = 0.01685393258426966
$ws.Range("E7").Value = 0.005617977528089887
$ws.Range("F7").Value = 0.06741573033707865
$ws.Range("J7").Value = 0.1123595505617977
$ws.Range("O7").Value = 0.03370786516853932
$ws.Range("Q7").Value = 0.1910112359550562
$ws.Range("R7").Value = 0.06741573033707865
$ws.Range("S7").Value = 0.3370786516853932
$ws.Range("B8").Value = 0.1113490364025696
$ws.Range("D8").Value = 0.03211991434689507
$ws.Range("E8").Value = 0.004282655246252677
$ws.Range("F8").Value = 0.05353319057815846
$ws.Range("J8").Value = 0.1027837259100642
$ws.Range("O8").Value = 0.008565310492505354
$ws.Range("Q8").Value = 0.2184154175588865
$ws.Range("R8").Value = 0.07494646680942184
$ws.Range("S8").Value = 0.3940042826552462
$ws.Range("B9").Value = 0.1375
$ws.Range("D9").Value = 0.0375
$ws.Range("F9").Value = 0.075
$ws.Range("J9").Value = 0.1125
$ws.Range("O9").Value = 0.00625
$ws.Range("Q9").Value = 0.18125
$ws.Range("R9").Value = 0.1
$ws.Range("S9").Value = 0.35
$ws.Range("B10").Value = 0.1147151898734177
$ws.Range("D10").Value = 0.0245253164556962
$ws.Range("F10").Value = 0.06408227848101265
$ws.Range("J10").Value = 0.1083860759493671
$ws.Range("O10").Value = 0.01503164556962025
$ws.Range("Q10").Value = 0.2602848101265823
$ws.Range("R10").Value = 0.07199367088607594
$ws.Range("S10").Value = 0.3409810126582278
$ws.Range("G11").Value = 0.1455223880597015
$ws.Range("J11").Value = 0.09328358208955224
$ws.Range("K11").Value = 0.1828358208955224
$ws.Range("L11").Value = 0.5783582089552238
$ws.Range("G12").Value = 0.7169811320754716
$ws.Range("J12").Value = 0.2138364779874214
$ws.Range("K12").Value = 0.01257861635220126
$ws.Range("L12").Value = 0.03144654088050314
$ws.Range("S12").Value = 0.02515723270440252
$ws.Range("F15").Value = 0.004484304932735426
$ws.Range("H15").Value = 0.1704035874439462
$ws.Range("I15").Value = 0.06278026905829596
$ws.Range("J15").Value = 0.3991031390134529
$ws.Range("K15").Value = 0.06278026905829596
$ws.Range("M15").Value = 0.03139013452914798
$ws.Range("O15").Value = 0.07623318385650224
$ws.Range("S15").Value = 0.1928251121076233
$ws.Range("F16").Value = 0.008968609865470852
$ws.Range("H16").Value = 0.2286995515695067
$ws.Range("I16").Value = 0.05381165919282511
$ws.Range("J16").Value = 0.3766816143497758
$ws.Range("K16").Value = 0.09417040358744394
$ws.Range("M16").Value = 0.02690582959641256
$ws.Range("N16").Value = 0.004484304932735426
$ws.Range("O16").Value = 0.06726457399103139
$ws.Range("S16").Value = 0.1390134529147982
$ws.Range("F17").Value = 0.01717557251908397
$ws.Range("H17").Value = 0.1965648854961832
$ws.Range("I17").Value = 0.07251908396946564
$ws.Range("J17").Value = 0.4427480916030535
$ws.Range("K17").Value = 0.09541984732824428
$ws.Range("M17").Value = 0.01526717557251908
$ws.Range("O17").Value = 0.06297709923664122
$ws.Range("S17").Value = 0.09732824427480916
$ws.Range("F18").Value = 0.02958579881656805
$ws.Range("H18").Value = 0.2544378698224852
$ws.Range("I18").Value = 0.07692307692307693
$ws.Range("J18").Value = 0.3609467455621302
$ws.Range("K18").Value = 0.1005917159763314
$ws.Range("M18").Value = 0.01775147928994083
$ws.Range("O18").Value = 0.07100591715976332
$ws.Range("S18").Value = 0.08875739644970414
$ws.Range("F19").Value = 0.01595744680851064
$ws.Range("H19").Value = 0.2083333333333333
$ws.Range("I19").Value = 0.07535460992907801
$ws.Range("J19").Value = 0.3953900709219858
$ws.Range("K19").Value = 0.09574468085106383
$ws.Range("M19").Value = 0.02304964539007092
$ws.Range("O19").Value = 0.07801418439716312
$ws.Range("S19").Value = 0.1081560283687943
